$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value would otherwise be auto-parsed as a number by Excel;
# force text format, assign, then restore default style so no stray numeric value/format remains.
$ws.Range("D2").Value = "61.836.36"
$ws.Range("E2").Value = "  +2.05%  "
$ws.Range("D3").Value = "2.413.79"
$ws.Range("E3").Value = "  +0.32%  "
$ws.Range("E4").Value = "  +0.65%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "569.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.08"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.89%  "
$ws.Range("E7").Value = "  -0.47%  "
$ws.Range("E8").Value = "  +0.69%  "
$ws.Range("D9").Value = "2.428.15"
$ws.Range("E9").Value = "  +1.62%  "
$ws.Range("E10").Value = "  +4.02%  "
$ws.Range("E11").Value = "  +0.63%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.23"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.20%  "
$ws.Range("E13").Value = "  +4.14%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.60"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.89%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000175"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.33%  "
$ws.Range("D17").Value = "61.744.41"
$ws.Range("E17").Value = "  +1.79%  "
$ws.Range("D18").Value = "2.428.28"
$ws.Range("E18").Value = "  +1.58%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.94"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.97%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.75"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.50%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "325.47"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.12%  "
$ws.Range("E22").Value = "  +1.75%  "
$ws.Range("E23").Value = "  +14.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.07"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.51%  "
$ws.Range("E25").Value = "  -0.19%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "65.20"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.77%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "623.65"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +13.56%  "
$ws.Range("E28").Value = "  +2.34%  "
$ws.Range("D29").Value = "0.0₃0964"
$ws.Range("E29").Value = "  +5.44%  "
$ws.Range("D30").Value = "2.521.37"
$ws.Range("E30").Value = "  -0.25%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.07"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.25%  "
$ws.Range("E32").Value = "  +9.00%  "
$ws.Range("E33").Value = "  +1.20%  "
$ws.Range("E34").Value = "  +2.23%  "
$ws.Range("E35").Value = "  +5.73%  "
$ws.Range("E36").Value = "  -0.59%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "153.25"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.50%  "
$ws.Range("E38").Value = "  +2.61%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.372"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.61%  "
$ws.Range("E40").Value = "  +5.70%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "18.45"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.60"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +11.58%  "
$ws.Range("E43").Value = "  +4.87%  "
$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "42.33"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.17%  "
$ws.Range("B45").Value = "USDe"
$ws.Range("C45").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.999"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.05%  "
$ws.Range("D46").Value = "0.0₆0285"
$ws.Range("E46").Value = "  -2.30%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "142.60"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.63%  "
$ws.Range("E48").Value = "  +2.07%  "
$ws.Range("E49").Value = "  +5.13%  "
$ws.Range("E50").Value = "  +2.70%  "
$ws.Range("E51").Value = "  +2.53%  "
